# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the "Conversión del día" text cell (A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$nl = [char]10
$conversionText = "Conversión del día 💰" + $nl + `
"✅ Dólar paralelo: 68" + $nl + `
$nl + `
"Binance" + $nl + `
"✅ 1000 Bs = 1.8 = 6547.94 pesos" + $nl + `
"✅ 6547.94 pesos = 1.79 = 940.55 Bs" + $nl + `
$nl + `
"Promedio competencia" + $nl + `
"✅ Tasa pesos: 20" + $nl + `
"✅ Tasa Bs: 20" + $nl + `
"✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $conversionText

# --- Sheet "tasas": update the rate cells N10, O10, N12, O12 ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 555.9
$wsTasas.Range("O10").Value = 3640
$wsTasas.Range("N12").Value = 3654.96
$wsTasas.Range("O12").Value = 525
